# Regenerate save_data to use K (column G) instead of Strike# values.
# Writes recalculated K values into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = 2
    16 = 2
    18 = 1
    19 = 1
    20 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
